$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("D1").Value = "Per Hour"
$ws.Range("E1").Value = "Earning $"

# Row 2 updates
$ws.Range("C2").Value = "SMU"
$ws.Range("D2").Value = 8
$ws.Range("E2").Value = 80

# Row 3 updates
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 100

# Row 4 new
$ws.Range("E4").Value = 180
